$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 30 (Fruta / Maracuyá, Primera
# quality), pushing the existing rows 30-86 down to 31-87.
$ws.Rows("30:30").Insert()

$ws.Range("A30").Value = 1
$ws.Range("B30").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C30").Value = "Arica y Parinacota"
$ws.Range("D30").Value = 44620
$ws.Range("E30").Value = 15
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100108
$ws.Range("H30").Value = "Tropicales y subtropicales"
$ws.Range("I30").Value = 100108003
$ws.Range("J30").Value = "Maracuyá"
$ws.Range("K30").Value = "Sin especificar"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 130
$ws.Range("N30").Value = 22000
$ws.Range("O30").Value = 23000
$ws.Range("P30").Value = 22500
$ws.Range("Q30").Value = "$/caja 20 kilos"
$ws.Range("R30").Value = "Región de Arica y Parinacota"
$ws.Range("S30").Value = 1125
$ws.Range("T30").Value = 20
